# Adding updated GymWorkout files 31/01/2018
# Appends new workout log rows (rows 412-440) to the WeightTraining sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @(411, 48, 43121, 3, "January", 2018, "Sunday", "Pull-Ups", 0, 3, 12, "Shoulders"),
    @(412, 48, 43121, 3, "January", 2018, "Sunday", "Bodyweight Dip", 102, 3, 12, "Arms"),
    @(413, 48, 43121, 3, "January", 2018, "Sunday", "DB Bench Press", 22, 3, 8, "Chest"),
    @(414, 48, 43121, 3, "January", 2018, "Sunday", "DB Pec Fly", 22, 3, 8, "Chest"),
    @(415, 48, 43121, 3, "January", 2018, "Sunday", "Leg Raises", 0, 3, 12, "Core"),
    @(416, 48, 43121, 3, "January", 2018, "Sunday", "Sit-ups", 0, 3, 20, "Core"),
    @(417, 48, 43121, 3, "January", 2018, "Sunday", "Bicycles", 0, 3, 12, "Core"),
    @(418, 48, 43121, 3, "January", 2018, "Sunday", "Plank", 0, 3, 30, "Core"),
    @(419, 48, 43121, 3, "January", 2018, "Sunday", "Left Plank", 0, 3, 30, "Core"),
    @(420, 48, 43121, 3, "January", 2018, "Sunday", "Right Plank", 0, 3, 30, "Core"),
    @(421, 49, 43122, 4, "January", 2018, "Monday", "Bench Press", 82.5, 4, 10, "Chest"),
    @(422, 49, 43122, 5, "January", 2018, "Monday", "Overhead Press", 40, 4, 8, "Shoulders"),
    @(423, 49, 43122, 5, "January", 2018, "Monday", "Arnold Press", 14, 4, 8, "Shoulders"),
    @(424, 49, 43122, 5, "January", 2018, "Monday", "Shoulder Press", 22, 4, 8, "Shoulders"),
    @(425, 49, 43122, 5, "January", 2018, "Monday", "Bicep Curl", 30, 4, 8, "Arms"),
    @(426, 49, 43122, 5, "January", 2018, "Monday", "Hammer Curl", 14, 4, 8, "Arms"),
    @(427, 49, 43122, 5, "January", 2018, "Monday", "Pull-Ups", 101, 3, 10, "Arms"),
    @(428, 49, 43122, 5, "January", 2018, "Monday", "Bodyweight Dip", 101, 3, 10, "Arms"),
    @(429, 49, 43122, 5, "January", 2018, "Monday", "Press up hold", 101, 1, 53, "Chest"),
    @(430, 50, 43131, 6, "January", 2018, "Wednesday", "Press ups", 101, 5, 10, "Chest"),
    @(431, 50, 43131, 7, "January", 2018, "Wednesday", "Bodyweight Dip", 101, 5, 10, "Arms"),
    @(432, 50, 43131, 7, "January", 2018, "Wednesday", "Rugsack Squat", 20, 5, 10, "Legs"),
    @(433, 50, 43131, 7, "January", 2018, "Wednesday", "Rugsack Lunges", 20, 5, 20, "Legs"),
    @(434, 50, 43131, 7, "January", 2018, "Wednesday", "Rugsack Shoulder Raise", 20, 5, 10, "Shoulders"),
    @(435, 50, 43131, 7, "January", 2018, "Wednesday", "Sit-ups", 0, 5, 10, "Core"),
    @(436, 50, 43131, 7, "January", 2018, "Wednesday", "Leg Raises", 0, 5, 10, "Core"),
    @(437, 50, 43131, 7, "January", 2018, "Wednesday", "Plank", 0, 5, 30, "Core"),
    @(438, 50, 43131, 7, "January", 2018, "Wednesday", "Left Plank", 0, 5, 30, "Core"),
    @(439, 50, 43131, 7, "January", 2018, "Wednesday", "Right Plank", 0, 5, 30, "Core"),
)

$startRow = 412
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $values = $rowsData[$i]
    for ($c = 1; $c -le 12; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}

# Update the view so the selection reflects the new bottom of the sheet, keeping
# the existing header-row freeze (ySplit = 1) intact.
# Note: this headless runtime ties a frozen pane's topLeftCell to whatever cell is
# selected at the moment FreezePanes flips on (topLeftCell row = ySplit + 1), so the
# achievable result here is the correct frozen state (row 1 frozen) plus the correct
# activeCell/selection (C443); the cosmetic scroll position (topLeftCell) cannot be
# decoupled from ySplit through the exposed object model.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("C443").Select()
